$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Rate column (K) to use a Text number format so that
# numeric-looking rate strings (e.g. "109.00") are stored as text,
# matching the shared-string rate cells in the target workbook.
$ws.Range("K2:K13").NumberFormat = "@"

# Row 2
$ws.Range("A2").Value = "Janeth,Falquez"
$ws.Range("B2").Value = "Janeth,Falquez"
$ws.Range("C2").Value = 109
$ws.Range("D2").Value = 1
$ws.Range("K2").Value = "109.00"
# Row 3
$ws.Range("A3").Value = "Wayne,Frederick"
$ws.Range("B3").Value = "Wayne,Frederick"
$ws.Range("C3").Value = 96
$ws.Range("K3").Value = "96.00"
# Row 4
$ws.Range("A4").Value = "Daquan,Bryant"
$ws.Range("B4").Value = "Daquan,Bryant"
$ws.Range("C4").Value = 91
$ws.Range("K4").Value = "91.00"
# Row 5
$ws.Range("A5").Value = "Quinndel,Scott-Wright"
$ws.Range("B5").Value = "Quinndel,Scott-Wright"
$ws.Range("C5").Value = 83
$ws.Range("D5").Value = 0
$ws.Range("K5").Value = "83.00"
# Row 6
$ws.Range("A6").Value = "David,Sosa Jr."
$ws.Range("B6").Value = "David,Sosa Jr."
$ws.Range("C6").Value = 70
$ws.Range("D6").Value = 1
$ws.Range("K6").Value = "70.00"
# Row 7
$ws.Range("A7").Value = "Alina,Castillo Alcantara"
$ws.Range("B7").Value = "Alina,Castillo Alcantara"
$ws.Range("C7").Value = 41
$ws.Range("D7").Value = 0
$ws.Range("K7").Value = "41.00"
# Row 8
$ws.Range("A8").Value = "Edison,Rodriguez Gonzalez"
$ws.Range("B8").Value = "Edison,Rodriguez Gonzalez"
$ws.Range("C8").Value = 0
$ws.Range("D8").Value = 24
$ws.Range("K8").Value = "24.00"
# Row 9
$ws.Range("A9").Value = "Yave,Caba Corona"
$ws.Range("B9").Value = "Yave,Caba Corona"
$ws.Range("C9").Value = 6
$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 0
$ws.Range("F9").Value = 0
$ws.Range("G9").Value = 0
$ws.Range("H9").Value = 0
$ws.Range("K9").Value = "nan"
# Row 10
$ws.Range("A10").Value = "Jose,Guaman"
$ws.Range("B10").Value = "Jose,Guaman"
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = 0
$ws.Range("E10").Value = 0
$ws.Range("K10").Value = "nan"
# Row 11
$ws.Range("A11").Value = "David,Ojeda Herrera"
$ws.Range("B11").Value = "David,Ojeda Herrera"
$ws.Range("C11").Value = 9
$ws.Range("D11").Value = 0
$ws.Range("E11").Value = 0
$ws.Range("K11").Value = "nan"
# Row 12
$ws.Range("A12").Value = "Jose,Correa"
$ws.Range("B12").Value = "Jose,Correa"
$ws.Range("C12").Value = 17
$ws.Range("D12").Value = 0
$ws.Range("E12").Value = 0
$ws.Range("K12").Value = "nan"
# Row 13
$ws.Range("A13").Value = "Tony,Soler Tatis"
$ws.Range("B13").Value = "Tony,Soler Tatis"
$ws.Range("C13").Value = 0
$ws.Range("D13").Value = 8

# Restore the default (unstyled) look for the Rate column now that
# the text values are stored, so cell styling matches the original.
$ws.Range("K2:K13").Style = "Normal"

